$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 0.984375
    3 = 0.96875
    4 = 0.96875
    5 = 0.984375
    6 = 0.96875
    7 = 0.953125
    8 = 0.953125
    9 = 1
    10 = 1
    11 = 1
    12 = 0.953125
    13 = 0.953125
    14 = 0.953125
    15 = 0.953125
    16 = 0.9375
    17 = 0.9375
    18 = 0.90625
    19 = 0.9375
    20 = 0.921875
    21 = 0.90625
    22 = 0.921875
    23 = 0.9375
    24 = 0.96875
    25 = 0.953125
    26 = 0.953125
    27 = 0.96875
    28 = 0.96875
    29 = 0.96875
    30 = 0.96875
    31 = 0.96875
    32 = 0.96875
    33 = 0.96875
    34 = 0.96875
    35 = 0.96875
    36 = 0.96875
    37 = 0.96875
    38 = 0.96875
    39 = 0.953125
    40 = 0.953125
    41 = 0.953125
    42 = 0.953125
    43 = 0.953125
    44 = 0.953125
    45 = 0.953125
    46 = 0.953125
    47 = 0.953125
    48 = 0.953125
    49 = 0.953125
    50 = 0.953125
    51 = 0.953125
    52 = 0.953125
    53 = 0.953125
    54 = 0.953125
    55 = 0.953125
    56 = 0.953125
    57 = 0.953125
    58 = 0.96875
    59 = 0.96875
    60 = 0.953125
    61 = 0.953125
    62 = 0.953125
    63 = 0.953125
    64 = 0.953125
    65 = 0.953125
    66 = 0.953125
    67 = 0.953125
    68 = 0.953125
    69 = 0.953125
    70 = 0.953125
    71 = 0.953125
    72 = 0.953125
    73 = 0.953125
    74 = 0.953125
    75 = 0.953125
    76 = 0.953125
    77 = 0.953125
    78 = 0.953125
    79 = 0.953125
    80 = 0.953125
    81 = 0.953125
    82 = 0.953125
    83 = 0.953125
    84 = 0.953125
    85 = 0.953125
    86 = 0.953125
    87 = 0.953125
    88 = 0.953125
    89 = 0.953125
    90 = 0.953125
    91 = 0.953125
    92 = 0.953125
    93 = 0.953125
    94 = 0.953125
    95 = 0.953125
    96 = 0.953125
    97 = 0.953125
    98 = 0.953125
    99 = 0.953125
    100 = 0.953125
    101 = 0.953125
    102 = 0.953125
    103 = 0.875
    104 = 0.953125
    105 = 0.921875
    106 = 0.953125
    107 = 0.921875
    108 = 0.90625
    109 = 0.921875
    110 = 0.984375
    111 = 0.921875
    112 = 0.953125
    113 = 0.984375
    114 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 2).Value = $values[$row]
}

